# Auto-generated Excel COM-interop script
# Applies the "ajustes carga e leiaute" edit: row 10 split into 4 rows
# (170512-Tocantins, 170520-RJ, 170522-Ceara, 170537-Maranhao),
# with the PLANO ESTADUAL block (AZ:BY) filled in per row and the
# leftover row 10 plano block cleared to "---".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Propagate the existing data-row style (border + centered + wrap) from
# row 10 down into the 3 freshly inserted rows before populating values.
$ws.Rows.Item(10).Copy()
$ws.Rows.Item(11).PasteSpecial(-4122)
$ws.Rows.Item(12).PasteSpecial(-4122)
$ws.Rows.Item(13).PasteSpecial(-4122)

# ---- Row 10 ----
$ws.Cells.Item(10, 1).Value = '170512'
$ws.Cells.Item(10, 2).Value = 'Departamento de Regulação, Avaliação e Controle de Sistemas DRAC'
$ws.Cells.Item(10, 3).Value = 'Programa de Redução de Filas e o Roteiro para Elaboração dos Planos Estaduais'
$ws.Cells.Item(10, 4).Value = 'PLANO ESTADUAL DE REDUÇÃO DE FILAS DE CIRURGIAS ELETIVAS'
$ws.Cells.Item(10, 5).Value = '---'
$ws.Cells.Item(10, 6).Value = 'Qualificado'
$ws.Cells.Item(10, 7).Value = 'Habilitação para custeio'
$ws.Cells.Item(10, 8).Value = '13.849.028/0001-40'
$ws.Cells.Item(10, 9).Value = 'FUNDO ESTADUAL DE SAUDE DO TOCANTINS'
$ws.Cells.Item(10, 10).Value = 'PALMAS'
$ws.Cells.Item(10, 11).Value = 'TO'
$ws.Cells.Item(10, 12).Value = 'Teste'
$ws.Cells.Item(10, 13).Value = 'Incompleta'
$ws.Cells.Item(10, 14).Value = '---'
$ws.Cells.Item(10, 15).Value = '---'
$ws.Cells.Item(10, 16).Value = '---'
$ws.Cells.Item(10, 17).Value = 'Regional'
$ws.Cells.Item(10, 18).Value = '---'
$ws.Cells.Item(10, 19).Value = 'ABREULANDIA; AGUIARNOPOLIS; ALIANCA DO TOCANTINS; ALMAS; ALVORADA; ANANAS; ANGICO; APARECIDA DO RIO NEGRO; ARAGOMINAS; ARAGUACEMA; ARAGUACU; ARAGUAINA; ARAGUANA; ARAGUATINS; ARAPOEMA; ARRAIAS; AUGUSTINOPOLIS; AURORA DO TOCANTINS; AXIXA DO TOCANTINS; BABACULANDIA; BANDEIRANTES DO TOCANTINS; BARRA DO OURO; BARROLANDIA; BERNARDO SAYAO; BOM JESUS DO TOCANTINS; BRASILANDIA DO TOCANTINS; BREJINHO DE NAZARE; BURITI DO TOCANTINS; CACHOEIRINHA; CAMPOS LINDOS; CARIRI DO TOCANTINS; CARMOLANDIA; CARRASCO BONITO; CASEARA; CENTENARIO; CHAPADA DE AREIA; CHAPADA DA NATIVIDADE; COLINAS DO TOCANTINS; COMBINADO; CONCEICAO DO TOCANTINS; COUTO MAGALHAES; CRISTALANDIA; CRIXAS DO TOCANTINS; DARCINOPOLIS; DIANOPOLIS; DIVINOPOLIS DO TOCANTINS; DOIS IRMAOS DO TOCANTINS; DUERE; ESPERANTINA; FATIMA; FIGUEIROPOLIS; FILADELFIA; FORMOSO DO ARAGUAIA; FORTALEZA DO TABOCAO; GOIANORTE; GOIATINS; GUARAI; GURUPI; IPUEIRAS; ITACAJA; ITAGUATINS; ITAPIRATINS; ITAPORA DO TOCANTINS; JAU DO TOCANTINS; JUARINA; LAGOA DA CONFUSAO; LAGOA DO TOCANTINS; LAJEADO; LAVANDEIRA; LIZARDA; LUZINOPOLIS; MARIANOPOLIS DO TOCANTINS; MATEIROS; MAURILANDIA DO TOCANTINS; MIRACEMA DO TOCANTINS; MIRANORTE; MONTE DO CARMO; MONTE SANTO DO TOCANTINS; PALMEIRAS DO TOCANTINS; MURICILANDIA; NATIVIDADE; NAZARE; NOVA OLINDA; NOVA ROSALANDIA; NOVO ACORDO; NOVO ALEGRE; NOVO JARDIM; OLIVEIRA DE FATIMA; PALMEIRANTE; PALMEIROPOLIS; PARAISO DO TOCANTINS; PARANA; PAU D''ARCO; PEDRO AFONSO; PEIXE; PEQUIZEIRO; COLMEIA; PINDORAMA DO TOCANTINS; PIRAQUE; PIUM; PONTE ALTA DO BOM JESUS; PONTE ALTA DO TOCANTINS; PORTO ALEGRE DO TOCANTINS; PORTO NACIONAL; PRAIA NORTE; PRESIDENTE KENNEDY; PUGMIL; RECURSOLANDIA; RIACHINHO; RIO DA CONCEICAO; RIO DOS BOIS; RIO SONO; SAMPAIO; SANDOLANDIA; SANTA FE DO ARAGUAIA; SANTA MARIA DO TOCANTINS; SANTA RITA DO TOCANTINS; SANTA ROSA DO TOCANTINS; SANTA TEREZA DO TOCANTINS; SANTA TEREZINHA DO TOCANTINS; SAO BENTO DO TOCANTINS; SAO FELIX DO TOCANTINS; SAO MIGUEL DO TOCANTINS; SAO SALVADOR DO TOCANTINS; SAO SEBASTIAO DO TOCANTINS; SAO VALERIO; SILVANOPOLIS; SITIO NOVO DO TOCANTINS; SUCUPIRA; TAGUATINGA; TAIPAS DO TOCANTINS; TALISMA; PALMAS; TOCANTINIA; TOCANTINOPOLIS; TUPIRAMA; TUPIRATINS; WANDERLANDIA; XAMBIOA'
$ws.Cells.Item(10, 20).Value = '170025; 170030; 170035; 170040; 170070; 170100; 170105; 170110; 170130; 170190; 170200; 170210; 170215; 170220; 170230; 170240; 170255; 170270; 170290; 170300; 170305; 170307; 170310; 170320; 170330; 170360; 170370; 170380; 170382; 170384; 170386; 170388; 170389; 170390; 170410; 170460; 170510; 170550; 170555; 170560; 170600; 170610; 170625; 170650; 170700; 170710; 170720; 170730; 170740; 170755; 170765; 170770; 170820; 170825; 170830; 170900; 170930; 170950; 170980; 171050; 171070; 171090; 171110; 171150; 171180; 171190; 171195; 171200; 171215; 171240; 171245; 171250; 171270; 171280; 171320; 171330; 171360; 171370; 171380; 171395; 171420; 171430; 171488; 171500; 171510; 171515; 171525; 171550; 171570; 171575; 171610; 171620; 171630; 171650; 171660; 171665; 171670; 171700; 171720; 171750; 171780; 171790; 171800; 171820; 171830; 171840; 171845; 171850; 171855; 171865; 171870; 171875; 171880; 171884; 171886; 171888; 171889; 171890; 171900; 172000; 172010; 172015; 172020; 172025; 172030; 172049; 172065; 172080; 172085; 172090; 172093; 172097; 172100; 172110; 172120; 172125; 172130; 172208; 172210'
$ws.Cells.Item(10, 21).Value = '1607363'
$ws.Cells.Item(10, 22).Value = '---'
$ws.Cells.Item(10, 23).Value = '---'
$ws.Cells.Item(10, 24).Value = '---'
$ws.Cells.Item(10, 25).Value = '---'
$ws.Cells.Item(10, 26).Value = '---'
$ws.Cells.Item(10, 27).Value = '---'
$ws.Cells.Item(10, 28).Value = '---'
$ws.Cells.Item(10, 29).Value = '---'
$ws.Cells.Item(10, 30).Value = '---'
$ws.Cells.Item(10, 31).Value = '---'
$ws.Cells.Item(10, 32).Value = '---'
$ws.Cells.Item(10, 33).Value = '---'
$ws.Cells.Item(10, 34).Value = '---'
$ws.Cells.Item(10, 35).Value = '---'
$ws.Cells.Item(10, 36).Value = '---'
$ws.Cells.Item(10, 37).Value = '---'
$ws.Cells.Item(10, 38).Value = '---'
$ws.Cells.Item(10, 39).Value = '---'
$ws.Cells.Item(10, 40).Value = '---'
$ws.Cells.Item(10, 41).Value = '---'
$ws.Cells.Item(10, 42).Value = '---'
$ws.Cells.Item(10, 43).Value = '---'
$ws.Cells.Item(10, 44).Value = 'JULIANA RIBEIRO PINTO'
$ws.Cells.Item(10, 45).Value = 'juenf86@gmail.com'
$ws.Cells.Item(10, 46).Value = '61 86119188'
$ws.Cells.Item(10, 47).Value = '24/02/2023'
$ws.Cells.Item(10, 48).Value = '25/05/2023'
$ws.Cells.Item(10, 49).Value = '---'
$ws.Cells.Item(10, 50).Value = '---'
$ws.Cells.Item(10, 51).Value = '---'
$ws.Cells.Item(10, 52).Value = '---'
$ws.Cells.Item(10, 53).Value = '---'
$ws.Cells.Item(10, 54).Value = '---'
$ws.Cells.Item(10, 55).Value = '---'
$ws.Cells.Item(10, 56).Value = '---'
$ws.Cells.Item(10, 57).Value = '---'
$ws.Cells.Item(10, 58).Value = '---'
$ws.Cells.Item(10, 59).Value = '---'
$ws.Cells.Item(10, 60).Value = '---'
$ws.Cells.Item(10, 61).Value = '---'
$ws.Cells.Item(10, 62).Value = '---'
$ws.Cells.Item(10, 63).Value = '---'
$ws.Cells.Item(10, 64).Value = '---'
$ws.Cells.Item(10, 65).Value = '---'
$ws.Cells.Item(10, 66).Value = '---'
$ws.Cells.Item(10, 67).Value = '---'
$ws.Cells.Item(10, 68).Value = '---'
$ws.Cells.Item(10, 69).Value = '---'
$ws.Cells.Item(10, 70).Value = '---'
$ws.Cells.Item(10, 71).Value = '---'
$ws.Cells.Item(10, 72).Value = '---'
$ws.Cells.Item(10, 73).Value = '---'
$ws.Cells.Item(10, 74).Value = '---'
$ws.Cells.Item(10, 75).Value = '---'
$ws.Cells.Item(10, 76).Value = '---'
$ws.Cells.Item(10, 77).Value = '---'

# ---- Row 11 ----
$ws.Cells.Item(11, 1).Value = '170520'
$ws.Cells.Item(11, 2).Value = 'Departamento de Regulação, Avaliação e Controle de Sistemas DRAC'
$ws.Cells.Item(11, 3).Value = 'Programa de Redução de Filas e o Roteiro para Elaboração dos Planos Estaduais'
$ws.Cells.Item(11, 4).Value = 'PLANO ESTADUAL DE REDUÇÃO DE FILAS DE CIRURGIAS ELETIVAS'
$ws.Cells.Item(11, 5).Value = '---'
$ws.Cells.Item(11, 6).Value = 'Qualificado'
$ws.Cells.Item(11, 7).Value = 'Habilitação para custeio'
$ws.Cells.Item(11, 8).Value = '35.949.791/0001-85'
$ws.Cells.Item(11, 9).Value = 'FUNDO ESTADUAL DE SAUDE FES'
$ws.Cells.Item(11, 10).Value = 'RIO DE JANEIRO'
$ws.Cells.Item(11, 11).Value = 'RJ'
$ws.Cells.Item(11, 12).Value = 'Plano teste.'
$ws.Cells.Item(11, 13).Value = 'Incompleta'
$ws.Cells.Item(11, 14).Value = '---'
$ws.Cells.Item(11, 15).Value = '---'
$ws.Cells.Item(11, 16).Value = '---'
$ws.Cells.Item(11, 17).Value = 'Regional'
$ws.Cells.Item(11, 18).Value = '---'
$ws.Cells.Item(11, 19).Value = 'BELFORD ROXO; DUQUE DE CAXIAS; ITAGUAI; JAPERI; MAGE; MESQUITA; NILOPOLIS; NOVA IGUACU; QUEIMADOS; RIO DE JANEIRO; SAO JOAO DE MERITI; SEROPEDICA'
$ws.Cells.Item(11, 20).Value = '330045; 330170; 330200; 330227; 330250; 330285; 330320; 330350; 330414; 330455; 330510; 330555'
$ws.Cells.Item(11, 21).Value = '10585667'
$ws.Cells.Item(11, 22).Value = '---'
$ws.Cells.Item(11, 23).Value = '---'
$ws.Cells.Item(11, 24).Value = '---'
$ws.Cells.Item(11, 25).Value = '---'
$ws.Cells.Item(11, 26).Value = '---'
$ws.Cells.Item(11, 27).Value = '---'
$ws.Cells.Item(11, 28).Value = '---'
$ws.Cells.Item(11, 29).Value = '---'
$ws.Cells.Item(11, 30).Value = '---'
$ws.Cells.Item(11, 31).Value = '---'
$ws.Cells.Item(11, 32).Value = '---'
$ws.Cells.Item(11, 33).Value = '---'
$ws.Cells.Item(11, 34).Value = '---'
$ws.Cells.Item(11, 35).Value = '---'
$ws.Cells.Item(11, 36).Value = '---'
$ws.Cells.Item(11, 37).Value = '---'
$ws.Cells.Item(11, 38).Value = '---'
$ws.Cells.Item(11, 39).Value = '---'
$ws.Cells.Item(11, 40).Value = '---'
$ws.Cells.Item(11, 41).Value = '---'
$ws.Cells.Item(11, 42).Value = '---'
$ws.Cells.Item(11, 43).Value = '---'
$ws.Cells.Item(11, 44).Value = 'LAURA VANESSA DE SOUZA ALBUQUERQUE'
$ws.Cells.Item(11, 45).Value = 'laura.albuquerque@saude.rj.gov.br'
$ws.Cells.Item(11, 46).Value = '21 23334031'
$ws.Cells.Item(11, 47).Value = '24/02/2023'
$ws.Cells.Item(11, 48).Value = '25/05/2023'
$ws.Cells.Item(11, 49).Value = '---'
$ws.Cells.Item(11, 50).Value = '---'
$ws.Cells.Item(11, 51).Value = '---'
$ws.Cells.Item(11, 52).Value = 'RJ'
$ws.Cells.Item(11, 53).Value = 'Laura Vanessa de Souza Albuquerque'
$ws.Cells.Item(11, 54).Value = 'Assistente Técnico'
$ws.Cells.Item(11, 55).Value = '(21) 2333-4031 
(21) 98561-9302 - Superintendente Marcelo Rodrigues
(21) 97015-6301 Laura'
$ws.Cells.Item(11, 56).Value = 'saecases@gmail.com
marcelo.rodrigues@saude.rj.gov.br
marcelo.rodrigues.castro@gmail.com
laura.albuquerque@saude.rj.gov.br'
$ws.Cells.Item(11, 57).Value = '---'
$ws.Cells.Item(11, 58).Value = '---'
$ws.Cells.Item(11, 59).Value = '---'
$ws.Cells.Item(11, 60).Value = '---'
$ws.Cells.Item(11, 61).Value = '---'
$ws.Cells.Item(11, 62).Value = '---'
$ws.Cells.Item(11, 63).Value = '---'
$ws.Cells.Item(11, 64).Value = '---'
$ws.Cells.Item(11, 65).Value = '---'
$ws.Cells.Item(11, 66).Value = '---'
$ws.Cells.Item(11, 67).Value = '---'
$ws.Cells.Item(11, 68).Value = '---'
$ws.Cells.Item(11, 69).Value = '---'
$ws.Cells.Item(11, 70).Value = '---'
$ws.Cells.Item(11, 71).Value = '---'
$ws.Cells.Item(11, 72).Value = '---'
$ws.Cells.Item(11, 73).Value = '---'
$ws.Cells.Item(11, 74).Value = '---'
$ws.Cells.Item(11, 75).Value = '---'
$ws.Cells.Item(11, 76).Value = '---'
$ws.Cells.Item(11, 77).Value = '---'

# ---- Row 12 ----
$ws.Cells.Item(12, 1).Value = '170522'
$ws.Cells.Item(12, 2).Value = 'Departamento de Regulação, Avaliação e Controle de Sistemas DRAC'
$ws.Cells.Item(12, 3).Value = 'Programa de Redução de Filas e o Roteiro para Elaboração dos Planos Estaduais'
$ws.Cells.Item(12, 4).Value = 'PLANO ESTADUAL DE REDUÇÃO DE FILAS DE CIRURGIAS ELETIVAS'
$ws.Cells.Item(12, 5).Value = '---'
$ws.Cells.Item(12, 6).Value = 'Qualificado'
$ws.Cells.Item(12, 7).Value = 'Habilitação para custeio'
$ws.Cells.Item(12, 8).Value = '74.031.865/0001-51'
$ws.Cells.Item(12, 9).Value = 'FUNDO ESTADUAL DE SAUDE'
$ws.Cells.Item(12, 10).Value = 'FORTALEZA'
$ws.Cells.Item(12, 11).Value = 'CE'
$ws.Cells.Item(12, 12).Value = 'Mutirão de cirurgias eletivas. Tem como objetivos: organizar e ampliar o acesso a cirurgias, exames e consultas na Atenção Especializada à 
Saúde, em especial àqueles com demanda reprimida identificada, aprimorar a governança da Rede de Atenção à Saúde com centralidade na garantia do acesso, gestão por resultados e financiamento estável; fomentar o monitoramento e a avaliação das ações e dos serviços de saúde, visando melhorar a qualidade da atenção especializada e ampliar o acesso à saúde; qualificar a contratualização com a rede complementar; mudar modelo de gestão e regulação das filas para a atenção especializada (regulação do acesso), visando a adequar a oferta de ações e serviços de saúde de acordo com as necessidades de saúde, estratificação de risco e necessidades assistenciais e fomentar a implementação de um novo modelo de custeio para a atenção ambulatorial especializada e para a realização de cirurgias eletivas.'
$ws.Cells.Item(12, 13).Value = 'Incompleta'
$ws.Cells.Item(12, 14).Value = '---'
$ws.Cells.Item(12, 15).Value = '---'
$ws.Cells.Item(12, 16).Value = '---'
$ws.Cells.Item(12, 17).Value = 'Regional'
$ws.Cells.Item(12, 18).Value = '---'
$ws.Cells.Item(12, 19).Value = 'ABAIARA; ACARAPE; ACARAU; ACOPIARA; AIUABA; ALCANTARAS; ALTANEIRA; ALTO SANTO; AMONTADA; ANTONINA DO NORTE; APUIARES; AQUIRAZ; ARACATI; ARACOIABA; ARARENDA; ARARIPE; ARATUBA; ARNEIROZ; ASSARE; AURORA; BAIXIO; BANABUIU; BARBALHA; BARREIRA; BARRO; BARROQUINHA; BATURITE; BEBERIBE; BELA CRUZ; BOA VIAGEM; BREJO SANTO; CAMOCIM; CAMPOS SALES; CANINDE; CAPISTRANO; CARIDADE; CARIRE; CARIRIACU; CARIUS; CARNAUBAL; CASCAVEL; CATARINA; CATUNDA; CAUCAIA; CEDRO; CHAVAL; CHORO; CHOROZINHO; COREAU; CRATEUS; CRATO; CROATA; CRUZ; DEPUTADO IRAPUAN PINHEIRO; ERERE; EUSEBIO; FARIAS BRITO; FORQUILHA; FORTALEZA; FORTIM; FRECHEIRINHA; GENERAL SAMPAIO; GRACA; GRANJA; GRANJEIRO; GROAIRAS; GUAIUBA; GUARACIABA DO NORTE; GUARAMIRANGA; HIDROLANDIA; HORIZONTE; IBARETAMA; IBIAPINA; IBICUITINGA; ICAPUI; ICO; IGUATU; INDEPENDENCIA; IPAPORANGA; IPAUMIRIM; IPU; IPUEIRAS; IRACEMA; IRAUCUBA; ITAICABA; ITAITINGA; ITAPAGE; ITAPIPOCA; ITAPIUNA; ITAREMA; ITATIRA; JAGUARETAMA; JAGUARIBARA; JAGUARIBE; JAGUARUANA; JARDIM; JATI; JIJOCA DE JERICOACOARA; JUAZEIRO DO NORTE; JUCAS; LAVRAS DA MANGABEIRA; LIMOEIRO DO NORTE; MADALENA; MARACANAU; MARANGUAPE; MARCO; MARTINOPOLE; MASSAPE; MAURITI; MERUOCA; MILAGRES; MILHA; MIRAIMA; MISSAO VELHA; MOMBACA; MONSENHOR TABOSA; MORADA NOVA; MORAUJO; MORRINHOS; MUCAMBO; MULUNGU; NOVA OLINDA; NOVA RUSSAS; NOVO ORIENTE; OCARA; OROS; PACAJUS; PACATUBA; PACOTI; PACUJA; PALHANO; PALMACIA; PARACURU; PARAIPABA; PARAMBU; PARAMOTI; PEDRA BRANCA; PENAFORTE; PENTECOSTE; PEREIRO; PINDORETAMA; PIQUET CARNEIRO; PIRES FERREIRA; PORANGA; PORTEIRAS; POTENGI; POTIRETAMA; QUITERIANOPOLIS; QUIXADA; QUIXELO; QUIXERAMOBIM; QUIXERE; REDENCAO; RERIUTABA; RUSSAS; SABOEIRO; SALITRE; SANTANA DO ACARAU; SANTANA DO CARIRI; SANTA QUITERIA; SAO BENEDITO; SAO GONCALO DO AMARANTE; SAO JOAO DO JAGUARIBE; SAO LUIS DO CURU; SENADOR POMPEU; SENADOR SA; SOBRAL; SOLONOPOLE; TABULEIRO DO NORTE; TAMBORIL; TARRAFAS; TAUA; TEJUCUOCA; TIANGUA; TRAIRI; TURURU; UBAJARA; UMARI; UMIRIM; URUBURETAMA; URUOCA; VARJOTA; VARZEA ALEGRE; VICOSA DO CEARA'
$ws.Cells.Item(12, 20).Value = '230010; 230015; 230020; 230030; 230040; 230050; 230060; 230070; 230075; 230080; 230090; 230100; 230110; 230120; 230125; 230130; 230140; 230150; 230160; 230170; 230180; 230185; 230190; 230195; 230200; 230205; 230210; 230220; 230230; 230240; 230250; 230260; 230270; 230280; 230290; 230300; 230310; 230320; 230330; 230340; 230350; 230360; 230365; 230370; 230380; 230390; 230393; 230395; 230400; 230410; 230420; 230423; 230425; 230426; 230427; 230428; 230430; 230435; 230440; 230445; 230450; 230460; 230465; 230470; 230480; 230490; 230495; 230500; 230510; 230520; 230523; 230526; 230530; 230533; 230535; 230540; 230550; 230560; 230565; 230570; 230580; 230590; 230600; 230610; 230620; 230625; 230630; 230640; 230650; 230655; 230660; 230670; 230680; 230690; 230700; 230710; 230720; 230725; 230730; 230740; 230750; 230760; 230763; 230765; 230770; 230780; 230790; 230800; 230810; 230820; 230830; 230835; 230837; 230840; 230850; 230860; 230870; 230880; 230890; 230900; 230910; 230920; 230930; 230940; 230945; 230950; 230960; 230970; 230980; 230990; 231000; 231010; 231020; 231025; 231030; 231040; 231050; 231060; 231070; 231080; 231085; 231090; 231095; 231100; 231110; 231120; 231123; 231126; 231130; 231135; 231140; 231150; 231160; 231170; 231180; 231190; 231195; 231200; 231210; 231220; 231230; 231240; 231250; 231260; 231270; 231280; 231290; 231300; 231310; 231320; 231325; 231330; 231335; 231340; 231350; 231355; 231360; 231370; 231375; 231380; 231390; 231395; 231400; 231410'
$ws.Cells.Item(12, 21).Value = '9240580'
$ws.Cells.Item(12, 22).Value = '---'
$ws.Cells.Item(12, 23).Value = '---'
$ws.Cells.Item(12, 24).Value = '---'
$ws.Cells.Item(12, 25).Value = '---'
$ws.Cells.Item(12, 26).Value = '---'
$ws.Cells.Item(12, 27).Value = '---'
$ws.Cells.Item(12, 28).Value = '---'
$ws.Cells.Item(12, 29).Value = '---'
$ws.Cells.Item(12, 30).Value = '25.991.043,34'
$ws.Cells.Item(12, 31).Value = '25.991.043,34'
$ws.Cells.Item(12, 32).Value = '---'
$ws.Cells.Item(12, 33).Value = '---'
$ws.Cells.Item(12, 34).Value = '---'
$ws.Cells.Item(12, 35).Value = '---'
$ws.Cells.Item(12, 36).Value = '---'
$ws.Cells.Item(12, 37).Value = '---'
$ws.Cells.Item(12, 38).Value = '---'
$ws.Cells.Item(12, 39).Value = '---'
$ws.Cells.Item(12, 40).Value = '---'
$ws.Cells.Item(12, 41).Value = '---'
$ws.Cells.Item(12, 42).Value = '---'
$ws.Cells.Item(12, 43).Value = '---'
$ws.Cells.Item(12, 44).Value = 'QUELVIA DA SILVA LIMA'
$ws.Cells.Item(12, 45).Value = 'quelviaadm@gmail.com'
$ws.Cells.Item(12, 46).Value = '85 99657881'
$ws.Cells.Item(12, 47).Value = '24/02/2023'
$ws.Cells.Item(12, 48).Value = '25/05/2023'
$ws.Cells.Item(12, 49).Value = '---'
$ws.Cells.Item(12, 50).Value = '---'
$ws.Cells.Item(12, 51).Value = '---'
$ws.Cells.Item(12, 52).Value = '---'
$ws.Cells.Item(12, 53).Value = '---'
$ws.Cells.Item(12, 54).Value = '---'
$ws.Cells.Item(12, 55).Value = '---'
$ws.Cells.Item(12, 56).Value = '---'
$ws.Cells.Item(12, 57).Value = '---'
$ws.Cells.Item(12, 58).Value = '---'
$ws.Cells.Item(12, 59).Value = '---'
$ws.Cells.Item(12, 60).Value = '---'
$ws.Cells.Item(12, 61).Value = '---'
$ws.Cells.Item(12, 62).Value = '---'
$ws.Cells.Item(12, 63).Value = '---'
$ws.Cells.Item(12, 64).Value = '---'
$ws.Cells.Item(12, 65).Value = '---'
$ws.Cells.Item(12, 66).Value = '---'
$ws.Cells.Item(12, 67).Value = '---'
$ws.Cells.Item(12, 68).Value = '---'
$ws.Cells.Item(12, 69).Value = '---'
$ws.Cells.Item(12, 70).Value = '---'
$ws.Cells.Item(12, 71).Value = '---'
$ws.Cells.Item(12, 72).Value = '---'
$ws.Cells.Item(12, 73).Value = '---'
$ws.Cells.Item(12, 74).Value = '---'
$ws.Cells.Item(12, 75).Value = '---'
$ws.Cells.Item(12, 76).Value = '---'
$ws.Cells.Item(12, 77).Value = '---'

# ---- Row 13 ----
$ws.Cells.Item(13, 1).Value = '170537'
$ws.Cells.Item(13, 2).Value = 'Departamento de Regulação, Avaliação e Controle de Sistemas DRAC'
$ws.Cells.Item(13, 3).Value = 'Programa de Redução de Filas e o Roteiro para Elaboração dos Planos Estaduais'
$ws.Cells.Item(13, 4).Value = 'PLANO ESTADUAL DE REDUÇÃO DE FILAS DE CIRURGIAS ELETIVAS'
$ws.Cells.Item(13, 5).Value = '---'
$ws.Cells.Item(13, 6).Value = 'Qualificado'
$ws.Cells.Item(13, 7).Value = 'Habilitação para custeio'
$ws.Cells.Item(13, 8).Value = '06.023.953/0001-51'
$ws.Cells.Item(13, 9).Value = 'ESTADO DO MARANHAO - FUNDO ESTADUAL DE SAUDE / FES'
$ws.Cells.Item(13, 10).Value = 'SAO LUIS'
$ws.Cells.Item(13, 11).Value = 'MA'
$ws.Cells.Item(13, 12).Value = '---'
$ws.Cells.Item(13, 13).Value = 'Incompleta'
$ws.Cells.Item(13, 14).Value = '---'
$ws.Cells.Item(13, 15).Value = '---'
$ws.Cells.Item(13, 16).Value = '---'
$ws.Cells.Item(13, 17).Value = 'Regional'
$ws.Cells.Item(13, 18).Value = '---'
$ws.Cells.Item(13, 19).Value = 'ACAILANDIA; AFONSO CUNHA; AGUA DOCE DO MARANHAO; ALCANTARA; ALDEIAS ALTAS; ALTAMIRA DO MARANHAO; ALTO ALEGRE DO MARANHAO; ALTO ALEGRE DO PINDARE; ALTO PARNAIBA; AMAPA DO MARANHAO; AMARANTE DO MARANHAO; ANAJATUBA; ANAPURUS; APICUM-ACU; ARAGUANA; ARAIOSES; ARAME; ARARI; AXIXA; BACABAL; BACABEIRA; BACURI; BACURITUBA; BALSAS; BARAO DE GRAJAU; BARRA DO CORDA; BARREIRINHAS; BELAGUA; BELA VISTA DO MARANHAO; BENEDITO LEITE; BEQUIMAO; BERNARDO DO MEARIM; BOA VISTA DO GURUPI; BOM JARDIM; BOM JESUS DAS SELVAS; BOM LUGAR; BREJO; BREJO DE AREIA; BURITI; BURITI BRAVO; BURITICUPU; BURITIRANA; CACHOEIRA GRANDE; CAJAPIO; CAJARI; CAMPESTRE DO MARANHAO; CANDIDO MENDES; CANTANHEDE; CAPINZAL DO NORTE; CAROLINA; CARUTAPERA; CAXIAS; CEDRAL; CENTRAL DO MARANHAO; CENTRO DO GUILHERME; CENTRO NOVO DO MARANHAO; CHAPADINHA; CIDELANDIA; CODO; COELHO NETO; COLINAS; CONCEICAO DO LAGO-ACU; COROATA; CURURUPU; DAVINOPOLIS; DOM PEDRO; DUQUE BACELAR; ESPERANTINOPOLIS; ESTREITO; FEIRA NOVA DO MARANHAO; FERNANDO FALCAO; FORMOSA DA SERRA NEGRA; FORTALEZA DOS NOGUEIRAS; FORTUNA; GODOFREDO VIANA; GONCALVES DIAS; GOVERNADOR ARCHER; GOVERNADOR EDISON LOBAO; GOVERNADOR EUGENIO BARROS; GOVERNADOR LUIZ ROCHA; GOVERNADOR NEWTON BELLO; GOVERNADOR NUNES FREIRE; GRACA ARANHA; GRAJAU; GUIMARAES; HUMBERTO DE CAMPOS; ICATU; IGARAPE DO MEIO; IGARAPE GRANDE; IMPERATRIZ; ITAIPAVA DO GRAJAU; ITAPECURU MIRIM; ITINGA DO MARANHAO; JATOBA; JENIPAPO DOS VIEIRAS; JOAO LISBOA; JOSELANDIA; JUNCO DO MARANHAO; LAGO DA PEDRA; LAGO DO JUNCO; LAGO VERDE; LAGOA DO MATO; LAGO DOS RODRIGUES; LAGOA GRANDE DO MARANHAO; LAJEADO NOVO; LIMA CAMPOS; LORETO; LUIS DOMINGUES; MAGALHAES DE ALMEIDA; MARACACUME; MARAJA DO SENA; MARANHAOZINHO; MATA ROMA; MATINHA; MATOES; MATOES DO NORTE; MILAGRES DO MARANHAO; MIRADOR; MIRANDA DO NORTE; MIRINZAL; MONCAO; MONTES ALTOS; MORROS; NINA RODRIGUES; NOVA COLINAS; NOVA IORQUE; NOVA OLINDA DO MARANHAO; OLHO D''AGUA DAS CUNHAS; OLINDA NOVA DO MARANHAO; PACO DO LUMIAR; PALMEIRANDIA; PARAIBANO; PARNARAMA; PASSAGEM FRANCA; PASTOS BONS; PAULINO NEVES; PAULO RAMOS; PEDREIRAS; PEDRO DO ROSARIO; PENALVA; PERI MIRIM; PERITORO; PINDARE-MIRIM; PINHEIRO; PIO XII; PIRAPEMAS; POCAO DE PEDRAS; PORTO FRANCO; PORTO RICO DO MARANHAO; PRESIDENTE DUTRA; PRESIDENTE JUSCELINO; PRESIDENTE MEDICI; PRESIDENTE SARNEY; PRESIDENTE VARGAS; PRIMEIRA CRUZ; RAPOSA; RIACHAO; RIBAMAR FIQUENE; ROSARIO; SAMBAIBA; SANTA FILOMENA DO MARANHAO; SANTA HELENA; SANTA INES; SANTA LUZIA; SANTA LUZIA DO PARUA; SANTA QUITERIA DO MARANHAO; SANTA RITA; SANTANA DO MARANHAO; SANTO AMARO DO MARANHAO; SANTO ANTONIO DOS LOPES; SAO BENEDITO DO RIO PRETO; SAO BENTO; SAO BERNARDO; SAO DOMINGOS DO AZEITAO; SAO DOMINGOS DO MARANHAO; SAO FELIX DE BALSAS; SAO FRANCISCO DO BREJAO; SAO FRANCISCO DO MARANHAO; SAO JOAO BATISTA; SAO JOAO DO CARU; SAO JOAO DO PARAISO; SAO JOAO DO SOTER; SAO JOAO DOS PATOS; SAO JOSE DE RIBAMAR; SAO JOSE DOS BASILIOS; SAO LUIS; SAO LUIS GONZAGA DO MARANHAO; SAO MATEUS DO MARANHAO; SAO PEDRO DA AGUA BRANCA; SAO PEDRO DOS CRENTES; SAO RAIMUNDO DAS MANGABEIRAS; SAO RAIMUNDO DO DOCA BEZERRA; SAO ROBERTO; SAO VICENTE FERRER; SATUBINHA; SENADOR ALEXANDRE COSTA; SENADOR LA ROCQUE; SERRANO DO MARANHAO; SITIO NOVO; SUCUPIRA DO NORTE; SUCUPIRA DO RIACHAO; TASSO FRAGOSO; TIMBIRAS; TIMON; TRIZIDELA DO VALE; TUFILANDIA; TUNTUM; TURIACU; TURILANDIA; TUTOIA; URBANO SANTOS; VARGEM GRANDE; VIANA; VILA NOVA DOS MARTIRIOS; VITORIA DO MEARIM; VITORINO FREIRE; ZE DOCA'
$ws.Cells.Item(13, 20).Value = '210005; 210010; 210015; 210020; 210030; 210040; 210043; 210047; 210050; 210055; 210060; 210070; 210080; 210083; 210087; 210090; 210095; 210100; 210110; 210120; 210125; 210130; 210135; 210140; 210150; 210160; 210170; 210173; 210177; 210180; 210190; 210193; 210197; 210200; 210203; 210207; 210210; 210215; 210220; 210230; 210232; 210235; 210237; 210240; 210250; 210255; 210260; 210270; 210275; 210280; 210290; 210300; 210310; 210312; 210315; 210317; 210320; 210325; 210330; 210340; 210350; 210355; 210360; 210370; 210375; 210380; 210390; 210400; 210405; 210407; 210408; 210409; 210410; 210420; 210430; 210440; 210450; 210455; 210460; 210462; 210465; 210467; 210470; 210480; 210490; 210500; 210510; 210515; 210520; 210530; 210535; 210540; 210542; 210545; 210547; 210550; 210560; 210565; 210570; 210580; 210590; 210592; 210594; 210596; 210598; 210600; 210610; 210620; 210630; 210632; 210635; 210637; 210640; 210650; 210660; 210663; 210667; 210670; 210675; 210680; 210690; 210700; 210710; 210720; 210725; 210730; 210735; 210740; 210745; 210750; 210760; 210770; 210780; 210790; 210800; 210805; 210810; 210820; 210825; 210830; 210840; 210845; 210850; 210860; 210870; 210880; 210890; 210900; 210905; 210910; 210920; 210923; 210927; 210930; 210940; 210945; 210950; 210955; 210960; 210970; 210975; 210980; 210990; 211000; 211003; 211010; 211020; 211023; 211027; 211030; 211040; 211050; 211060; 211065; 211070; 211080; 211085; 211090; 211100; 211102; 211105; 211107; 211110; 211120; 211125; 211130; 211140; 211150; 211153; 211157; 211160; 211163; 211167; 211170; 211172; 211174; 211176; 211178; 211180; 211190; 211195; 211200; 211210; 211220; 211223; 211227; 211230; 211240; 211245; 211250; 211260; 211270; 211280; 211285; 211290; 211300; 211400'
$ws.Cells.Item(13, 21).Value = '7153262'
$ws.Cells.Item(13, 22).Value = '---'
$ws.Cells.Item(13, 23).Value = '---'
$ws.Cells.Item(13, 24).Value = '---'
$ws.Cells.Item(13, 25).Value = '---'
$ws.Cells.Item(13, 26).Value = '---'
$ws.Cells.Item(13, 27).Value = '---'
$ws.Cells.Item(13, 28).Value = '---'
$ws.Cells.Item(13, 29).Value = '---'
$ws.Cells.Item(13, 30).Value = '20.120.029,55'
$ws.Cells.Item(13, 31).Value = '20.120.029,55'
$ws.Cells.Item(13, 32).Value = '---'
$ws.Cells.Item(13, 33).Value = '---'
$ws.Cells.Item(13, 34).Value = '---'
$ws.Cells.Item(13, 35).Value = '---'
$ws.Cells.Item(13, 36).Value = '---'
$ws.Cells.Item(13, 37).Value = '---'
$ws.Cells.Item(13, 38).Value = '---'
$ws.Cells.Item(13, 39).Value = '---'
$ws.Cells.Item(13, 40).Value = '---'
$ws.Cells.Item(13, 41).Value = '---'
$ws.Cells.Item(13, 42).Value = '---'
$ws.Cells.Item(13, 43).Value = '---'
$ws.Cells.Item(13, 44).Value = 'LUCIANA AMORIM TOMICH NETTO GUTERRES SOARES'
$ws.Cells.Item(13, 45).Value = 'lucianatomich@yahoo.com.br'
$ws.Cells.Item(13, 46).Value = '98 99143176'
$ws.Cells.Item(13, 47).Value = '27/02/2023'
$ws.Cells.Item(13, 48).Value = '28/05/2023'
$ws.Cells.Item(13, 49).Value = '---'
$ws.Cells.Item(13, 50).Value = '---'
$ws.Cells.Item(13, 51).Value = '---'
$ws.Cells.Item(13, 52).Value = 'MA'
$ws.Cells.Item(13, 53).Value = 'Marina Nascimento Sousa'
$ws.Cells.Item(13, 54).Value = 'Superintendência de Avaliação e Controle do Sistema de Saúde/SES'
$ws.Cells.Item(13, 55).Value = '(98) 98127-8774'
$ws.Cells.Item(13, 56).Value = 'marinasousa20@gmail.com'
$ws.Cells.Item(13, 57).Value = '---'
$ws.Cells.Item(13, 58).Value = '---'
$ws.Cells.Item(13, 59).Value = 'Há uma demanda reprimida para procedimentos cirúrgicos eletivos de media e alta complexidade, potencializada pela Pandemia Covid 19. A Secretaria de Estado da Saúde do Maranhão em conjunto com os municípios descentralizou os serviços de saúde ambulatoriais e hospitalares para as 19 Regiões de Saúde, abrangendo os 217 municípios maranhenses. Desta feita, os procedimentos cirúrgicos eletivos serão realizados nos estabelecimentos de saúde de gestão estadual e municipal'
$ws.Cells.Item(13, 60).Value = '---'
$ws.Cells.Item(13, 61).Value = 'Sim'
$ws.Cells.Item(13, 62).Value = 'Sim'
$ws.Cells.Item(13, 63).Value = 'Sim'
$ws.Cells.Item(13, 64).Value = 'Sim'
$ws.Cells.Item(13, 65).Value = 'Sim'
$ws.Cells.Item(13, 66).Value = 'Sim'
$ws.Cells.Item(13, 67).Value = 'Sim'
$ws.Cells.Item(13, 68).Value = 'Sim'
$ws.Cells.Item(13, 69).Value = 'Não'
$ws.Cells.Item(13, 70).Value = 'Não'
$ws.Cells.Item(13, 71).Value = 'Não'
$ws.Cells.Item(13, 72).Value = 'Para ampliação dos serviços será necessário implementar os serviços nos estabelecimentos de saude adotando as seguintes medidas:
* Aumento da capacidade de atendiemnto dos serviços de saude no ambito ambulatorial e hospitalar;
* Melhor utilização da capacidade dos Serviços implementando o sistema de produtividade;
* Ampliação de turno de trabalho para os procedimentos a serem realizados;
* Ampliação das equipes multiprofissionais de saúde para a realização dos procedimentos cirúrgicos.'
$ws.Cells.Item(13, 73).Value = '---'
$ws.Cells.Item(13, 74).Value = 'A gestão da fila cirúrgica será realizada pela gestão estadual com sistema de compartilhamento com os municípios, visando o acompanhamento das metas programadas e realizadas por cada estabelecimento de saúde e cronograma de execução estabelecido no Plano Nacional de Redução de Filas'
$ws.Cells.Item(13, 75).Value = '---'
$ws.Cells.Item(13, 76).Value = '---'
$ws.Cells.Item(13, 77).Value = '---'

# Match the workbook author selection left on the sheet after editing.
$ws.Range("BY13").Select()
